$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "46.366.43"
$ws.Range("E2").Value = "  -1.15%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.614.96"
$ws.Range("E3").Value = "  +0.91%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.01"
$ws.Range("E5").Value = "  +0.44%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "100.11"
$ws.Range("E6").Value = "  -2.49%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.602"
$ws.Range("E7").Value = "  -0.37%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  +1.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.35"
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0845"
$ws.Range("E11").Value = "  +1.09%  "
$ws.Range("B12").Value = "OKB"
$ws.Range("C12").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.32"
$ws.Range("E12").Value = "  -1.31%  "
$ws.Range("E13").Value = "  +2.20%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.014.54"
$ws.Range("E14").Value = "  +0.88%  "
$ws.Range("E15").Value = "  +0.66%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.617.25"
$ws.Range("E16").Value = "  +0.11%  "
$ws.Range("E17").Value = "  +2.76%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "15.02"
$ws.Range("E18").Value = "  -0.81%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "46.525.52"
$ws.Range("E19").Value = "  -1.23%  "
$ws.Range("E20").Value = "  +1.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.03"
$ws.Range("E21").Value = "  -5.83%  "
$ws.Range("E22").Value = "  +2.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "71.94"
$ws.Range("E23").Value = "  +2.44%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "275.58"
$ws.Range("E24").Value = "  +7.55%  "
$ws.Range("E25").Value = "  +1.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.21"
$ws.Range("E26").Value = "  +4.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "29.59"
$ws.Range("E27").Value = "  +13.20%  "
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.02"
$ws.Range("E29").Value = "  -1.64%  "
$ws.Range("E30").Value = "  +1.52%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "38.71"
$ws.Range("E31").Value = "  -6.56%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.22"
$ws.Range("E32").Value = "  -2.76%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.44"
$ws.Range("E33").Value = "  +8.14%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.64"
$ws.Range("E34").Value = "  -4.28%  "
$ws.Range("E35").Value = "  +0.90%  "
$ws.Range("E36").Value = "  -2.37%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0840"
$ws.Range("E37").Value = "  -0.54%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "152.30"
$ws.Range("E38").Value = "  +1.73%  "
$ws.Range("E39").Value = "  -0.04%  "
$ws.Range("E40").Value = "  +1.64%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "23.87"
$ws.Range("E41").Value = "  +30.18%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "15.99"
$ws.Range("E42").Value = "  -2.72%  "
$ws.Range("E43").Value = "  +0.36%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.63"
$ws.Range("E44").Value = "  +0.86%  "
$ws.Range("E45").Value = "  -4.86%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.137.45"
$ws.Range("E46").Value = "  +5.58%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.998"
$ws.Range("E47").Value = "  -0.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "94.63"
$ws.Range("E48").Value = "  +1.37%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.53"
$ws.Range("E49").Value = "  +7.24%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "109.77"
$ws.Range("E50").Value = "  +1.13%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.78"
$ws.Range("E51").Value = "  -5.07%  "
